$d = $word.ActiveDocument

$d.Content.Find.Execute("207×8=1656", $true, $false, $false, $false, $false, $true, 1, $false, "448×8=3584", 2) | Out-Null
$d.Content.Find.Execute("489×5=2445", $true, $false, $false, $false, $false, $true, 1, $false, "556×9=5004", 2) | Out-Null
$d.Content.Find.Execute("369×8=2952", $true, $false, $false, $false, $false, $true, 1, $false, "799×7=5593", 2) | Out-Null
$d.Content.Find.Execute("487×2=974", $true, $false, $false, $false, $false, $true, 1, $false, "584×3=1752", 2) | Out-Null
$d.Content.Find.Execute("584×5=2920", $true, $false, $false, $false, $false, $true, 1, $false, "376×3=1128", 2) | Out-Null
$d.Content.Find.Execute("177×8=1416", $true, $false, $false, $false, $false, $true, 1, $false, "223×9=2007", 2) | Out-Null
$d.Content.Find.Execute("309×9=2781", $true, $false, $false, $false, $false, $true, 1, $false, "281×3=843", 2) | Out-Null
$d.Content.Find.Execute("607×6=3642", $true, $false, $false, $false, $false, $true, 1, $false, "945×2=1890", 2) | Out-Null
$d.Content.Find.Execute("625×5=3125", $true, $false, $false, $false, $false, $true, 1, $false, "450×4=1800", 2) | Out-Null
$d.Content.Find.Execute("636×3=1908", $true, $false, $false, $false, $false, $true, 1, $false, "733×7=5131", 2) | Out-Null
$d.Content.Find.Execute("656×6=3936", $true, $false, $false, $false, $false, $true, 1, $false, "496×4=1984", 2) | Out-Null
$d.Content.Find.Execute("305×4=1220", $true, $false, $false, $false, $false, $true, 1, $false, "102×5=510", 2) | Out-Null
$d.Content.Find.Execute("266×7=1862", $true, $false, $false, $false, $false, $true, 1, $false, "772×2=1544", 2) | Out-Null
$d.Content.Find.Execute("616×4=2464", $true, $false, $false, $false, $false, $true, 1, $false, "452×8=3616", 2) | Out-Null
$d.Content.Find.Execute("219×2=438", $true, $false, $false, $false, $false, $true, 1, $false, "164×3=492", 2) | Out-Null
$d.Content.Find.Execute("785×8=6280", $true, $false, $false, $false, $false, $true, 1, $false, "557×4=2228", 2) | Out-Null
$d.Content.Find.Execute("117×5=585", $true, $false, $false, $false, $false, $true, 1, $false, "628×3=1884", 2) | Out-Null
$d.Content.Find.Execute("266×5=1330", $true, $false, $false, $false, $false, $true, 1, $false, "697×9=6273", 2) | Out-Null
$d.Content.Find.Execute("838×8=6704", $true, $false, $false, $false, $false, $true, 1, $false, "433×6=2598", 2) | Out-Null
$d.Content.Find.Execute("306×6=1836", $true, $false, $false, $false, $false, $true, 1, $false, "207×9=1863", 2) | Out-Null
$d.Content.Find.Execute("968×5=4840", $true, $false, $false, $false, $false, $true, 1, $false, "485×3=1455", 2) | Out-Null
$d.Content.Find.Execute("388×7=2716", $true, $false, $false, $false, $false, $true, 1, $false, "519×2=1038", 2) | Out-Null
$d.Content.Find.Execute("824×2=1648", $true, $false, $false, $false, $false, $true, 1, $false, "982×8=7856", 2) | Out-Null
$d.Content.Find.Execute("114×3=342", $true, $false, $false, $false, $false, $true, 1, $false, "298×3=894", 2) | Out-Null
$d.Content.Find.Execute("622×4=2488", $true, $false, $false, $false, $false, $true, 1, $false, "369×6=2214", 2) | Out-Null

Write-Host "Replacements complete"
